$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update UnitType comment text, clear old B1 comment
$ws.Range("A1").Value = "UnitType(EMS_UnitType) - 1 : BasePlayer, 2 : Item 3 : Furniture "
$ws.Range("B1").Value = $null

# New row 2 for the ChildTableDataType comment
$ws.Range("A2").Value = "ChildTableDataType(EMS_TableDataType) - 10 : ItemData 12 : Storage"
$ws.Range("A2").Style = $ws.Range("A1").Style

# Row 3 headers: Base_Path -> UnitType, UnitType -> ChildTableDataType, clear Grid columns
$ws.Range("B3").Value = "UnitType"
$ws.Range("C3").Value = "ChildTableDataType"
$ws.Range("D3").Value = $null
$ws.Range("E3").Value = $null

# Row 4 values
$ws.Range("C4").Value = -1
$ws.Range("D4").Value = $null
$ws.Range("E4").Value = $null

# Row 5: fill in the rest of the record (previously only C5 had a value)
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2

# New row 6: additional record
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 12

$ws.Range("G8").Select()
